$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 222
$ws.Range("F3").Value = 527
$ws.Range("G3").Value = "不可售"
$ws.Range("G4").Value = 65
$ws.Range("F5").Value = 131
$ws.Range("G5").Value = "不可售"
$ws.Range("F7").Value = 66
$ws.Range("F8").Value = 7116
$ws.Range("F9").Value = 259
$ws.Range("F11").Value = 3582
$ws.Range("F12").Value = 302
$ws.Range("F13").Value = 522
$ws.Range("F14").Value = 265
$ws.Range("F15").Value = 592
$ws.Range("F16").Value = 75

$ws = $wb.Worksheets.Item("演出")
$ws.Range("B3").Value = "2024-08-02"
$ws.Range("C3").Value = "合肥·新西兰·治愈系民谣歌手Luke Thompson2024中国巡演 KEEP ROLLING ON "
$ws.Range("D3").Value = "宁国路罍街二期15号楼安徽原创音乐基地3楼 合肥ON THE WAY LiveHouse"
$ws.Range("E3").Value = "2024.08.02 20:00-08.02 21:30"
$ws.Range("F3").Value = 4
$ws.Range("G3").Value = 180
$ws.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=88824"
$ws.Range("I3").Value = "//i2.hdslb.com/bfs/openplatform/202407/FKItJRNl1719803666645.jpeg"
$ws.Range("B4").Value = "2024-08-03"
$ws.Range("C4").Value = "合肥·首届包河留声机音乐节—《菊次郎的夏天》久石让钢琴曲梦幻之旅演奏会"
$ws.Range("D4").Value = "徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院"
$ws.Range("E4").Value = "2024.08.03 19:30-08.03 21:00"
$ws.Range("F4").Value = 49
$ws.Range("G4").Value = 80
$ws.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=83556"
$ws.Range("I4").Value = "//i1.hdslb.com/bfs/openplatform/202403/4nwOTVDu1711695345941.jpeg"
$ws.Rows.Item(5).Delete()

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("C3").Value = "庐江·夏日游嘉年华"
$ws.Range("D3").Value = "白山路东150米 庐江体育馆"
$ws.Range("E3").Value = "2024.07.27 09:00-07.28 17:00"
$ws.Range("F3").Value = 222
$ws.Range("G3").Value = 70
$ws.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=87569"
$ws.Range("I3").Value = "//i2.hdslb.com/bfs/openplatform/202406/5tB3RWrN1718243791381.jpeg"
$ws.Range("B4").Value = "2024-07-28"
$ws.Range("C4").Value = "合肥·咒术回战only"
$ws.Range("D4").Value = "清河路19号 依立腾工业园区"
$ws.Range("E4").Value = "2024.07.28 09:30-07.28 17:30"
$ws.Range("F4").Value = 527
$ws.Range("G4").Value = "不可售"
$ws.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=86520"
$ws.Range("I4").Value = "//i2.hdslb.com/bfs/openplatform/202405/cLCM0a1e1716952386781.png"
$ws.Range("C5").Value = "合肥·第二届TH动漫游戏嘉年华"
$ws.Range("D5").Value = "田埠西路199号 吉祥如意宴会楼蜀山店"
$ws.Range("E5").Value = "2024.07.28 09:30-07.28 17:00"
$ws.Range("F5").Value = 119
$ws.Range("G5").Value = 65
$ws.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=87447"
$ws.Range("I5").Value = "//i0.hdslb.com/bfs/openplatform/202406/jHqfdzLQ1718091324240.png"
$ws.Range("C6").Value = "合肥·首届进击的巨人ONLY漫展"
$ws.Range("D6").Value = "胜利路198号 合肥元一希尔顿酒店"
$ws.Range("E6").Value = "2024.07.28 09:30-07.28 16:30"
$ws.Range("F6").Value = 131
$ws.Range("G6").Value = "不可售"
$ws.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=88965"
$ws.Range("I6").Value = "//i1.hdslb.com/bfs/openplatform/202406/q9ANU7gh1718880973689.jpeg"
$ws.Range("B7").Value = "2024-07-30"
$ws.Range("C7").Value = "巢湖·元气动漫游戏嘉年华"
$ws.Range("D7").Value = "团结东路7号 巢湖宾馆"
$ws.Range("E7").Value = "2024.07.30 10:00-07.30 17:00"
$ws.Range("F7").Value = 53
$ws.Range("G7").Value = 45
$ws.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=88193"
$ws.Range("I7").Value = "//i2.hdslb.com/bfs/openplatform/202406/3VBeQfqQ1719318873395.jpeg"
$ws.Range("B8").Value = "2024-08-01"
$ws.Range("C8").Value = "合肥·ACGN夏日游园会预热场"
$ws.Range("D8").Value = "五里墩街道长江西路与金牛路交叉口向北300米 水善汇都市微度假"
$ws.Range("E8").Value = "2024.08.01 09:30-08.02 18:00"
$ws.Range("F8").Value = 66
$ws.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=89914"
$ws.Range("I8").Value = "//i2.hdslb.com/bfs/openplatform/202407/rfRl5Bbj1721815713827.jpeg"
$ws.Range("B9").Value = "2024-08-02"
$ws.Range("C9").Value = "合肥·新西兰·治愈系民谣歌手Luke Thompson2024中国巡演 KEEP ROLLING ON "
$ws.Range("D9").Value = "宁国路罍街二期15号楼安徽原创音乐基地3楼 合肥ON THE WAY LiveHouse"
$ws.Range("E9").Value = "2024.08.02 20:00-08.02 21:30"
$ws.Range("F9").Value = 4
$ws.Range("G9").Value = 180
$ws.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=88824"
$ws.Range("I9").Value = "//i2.hdslb.com/bfs/openplatform/202407/FKItJRNl1719803666645.jpeg"
$ws.Range("B10").Value = "2024-08-03"
$ws.Range("C10").Value = "合肥·第七届环形宇宙动漫游戏嘉年华"
$ws.Range("D10").Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws.Range("E10").Value = "2024.08.03 09:30-08.04 17:00"
$ws.Range("F10").Value = 7116
$ws.Range("G10").Value = 44.1
$ws.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=84767"
$ws.Range("I10").Value = "//i2.hdslb.com/bfs/openplatform/202404/nBGuQecO1713856894035.jpeg"
$ws.Range("C11").Value = "合肥·首届包河留声机音乐节—《菊次郎的夏天》久石让钢琴曲梦幻之旅演奏会"
$ws.Range("D11").Value = "徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院"
$ws.Range("E11").Value = "2024.08.03 19:30-08.03 21:00"
$ws.Range("F11").Value = 49
$ws.Range("G11").Value = 80
$ws.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=83556"
$ws.Range("I11").Value = "//i1.hdslb.com/bfs/openplatform/202403/4nwOTVDu1711695345941.jpeg"
$ws.Range("B12").Value = "2024-08-10"
$ws.Range("C12").Value = "合肥·排球少年only之夏日招新季"
$ws.Range("D12").Value = "广德路与长江东路交口往北200米文一时埠里文旅街区 巅峰篮球公园"
$ws.Range("E12").Value = "2024.08.10 10:00-08.10 17:00"
$ws.Range("F12").Value = 259
$ws.Range("G12").Value = 70
$ws.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=88281"
$ws.Range("I12").Value = "//i0.hdslb.com/bfs/openplatform/202406/qjd7yzXE1719556597555.jpeg"
$ws.Range("C13").Value = "合肥·比翼连枝国乙&代号鸢only"
$ws.Range("D13").Value = "长江东大街与东二环路交叉口向南300米东方摩域商业广场三楼 格律诗婚礼艺术中心(筑梦店)"
$ws.Range("E13").Value = "2024.08.10 09:00-08.10 22:00"
$ws.Range("F13").Value = 409
$ws.Range("G13").Value = 65
$ws.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=88421"
$ws.Range("I13").Value = "//i0.hdslb.com/bfs/openplatform/202407/RHiXT98J1721199172046.jpeg"
$ws.Range("B14").Value = "2024-08-17"
$ws.Range("C14").Value = "合肥·第八届环形宇宙动漫游戏嘉年华Plus"
$ws.Range("D14").Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws.Range("E14").Value = "2024.08.17 09:30-08.18 17:00"
$ws.Range("F14").Value = 3582
$ws.Range("G14").Value = 69
$ws.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=88650"
$ws.Range("I14").Value = "//i2.hdslb.com/bfs/openplatform/202407/4I7mduRV1720071650216.jpeg"
$ws.Range("C15").Value = "合肥·第八届环形宇宙动漫游戏嘉年华Plus~水千丞周边预约票"
$ws.Range("E15").Value = "2024.08.17 09:30-08.17 17:00"
$ws.Range("F15").Value = 302
$ws.Range("G15").Value = 0.1
$ws.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=89420"
$ws.Range("I15").Value = "//i0.hdslb.com/bfs/openplatform/202407/hsiXAged1721203655434.jpeg"
$ws.Range("C16").Value = "合肥·第八届环形宇宙动漫游戏嘉年华Plus~水千丞签售预约票"
$ws.Range("F16").Value = 522
$ws.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=89421"
$ws.Range("I16").Value = "//i2.hdslb.com/bfs/openplatform/202407/r8wJqvVO1721202573195.jpeg"
$ws.Range("C17").Value = "合肥·银魂主题派对only2.0"
$ws.Range("D17").Value = "长江东路1137号圣大国际商贸中心2-301室 梦田音乐LiveHouse(合肥店)"
$ws.Range("E17").Value = "2024.08.17 13:00-08.17 18:00"
$ws.Range("F17").Value = 265
$ws.Range("G17").Value = 128
$ws.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=87173"
$ws.Range("I17").Value = "//i2.hdslb.com/bfs/openplatform/202406/aSc8SoTl1718078234193.png"
$ws.Range("B18").Value = "2024-08-18"
$ws.Range("C18").Value = "合肥·SSS第五人格only"
$ws.Range("D18").Value = "桐城路127号合作经济广场3号楼23层 赤阑桥艺术空间"
$ws.Range("E18").Value = "2024.08.18 09:00-08.18 17:00"
$ws.Range("F18").Value = 592
$ws.Range("G18").Value = 68
$ws.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=88430"
$ws.Range("I18").Value = "//i0.hdslb.com/bfs/openplatform/202406/a0qh8I1h1719660853555.png"
$ws.Range("B19").Value = "2024-09-07"
$ws.Range("C19").Value = "合肥·国乙only宇宙心动（含夜场）"
$ws.Range("D19").Value = "文忠路1865号 赫拉诺言艺术中心"
$ws.Range("E19").Value = "2024.09.07 10:00-09.07 21:00"
$ws.Range("F19").Value = 75
$ws.Range("G19").Value = 48
$ws.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=89803"
$ws.Range("I19").Value = "//i1.hdslb.com/bfs/openplatform/202407/w5hQDj821721564303601.jpeg"
$ws.Rows.Item(20).Delete()
